# Week 5 day 1 update: add Discoverer / Year of Discovery / Composition
# columns, shorten the Feature column text, and append the Pluto row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells G1:I1 (styled like the existing F1 header: bold,
#    thin border all around, centered horizontally, top vertically).
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Discoverer"
$ws.Range("H1").Value = "Year of Discovery"
$ws.Range("I1").Value = "Composition"

$headerRng = $ws.Range("G1:I1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1         # xlContinuous

# ---------------------------------------------------------------------------
# 2. Shorten the existing "Feature" text in column F, and fill in the new
#    Discoverer / Year of Discovery / Composition columns for rows 2-9.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Eccentric orbit"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "N/A"
$ws.Range("I2").Value = "Iron"

$ws.Range("F3").Value = "runaway greenhouse"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "N/A"
$ws.Range("I3").Value = "Carbon dioxide"

$ws.Range("F4").Value = "only planet known to support life"
$ws.Range("G4").Value = "N/A"
$ws.Range("H4").Value = "N/A"
$ws.Range("I4").Value = "Nitrogen and oxygen"

$ws.Range("F5").Value = "largest volcano"
$ws.Range("G5").Value = "N/A"
$ws.Range("H5").Value = "N/A"
$ws.Range("I5").Value = "Iron oxide"

$ws.Range("F6").Value = "Red Spot is a massive storm"
$ws.Range("G6").Value = "N/A"
$ws.Range("H6").Value = "N/A"
$ws.Range("I6").Value = "Mostly hydrogen and helium"

$ws.Range("F7").Value = "Saturn's rings are the most extensive"
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = "N/A"
$ws.Range("I7").Value = "Mostly hydrogen and helium"

$ws.Range("F8").Value = "Uranus rotates on its side"
$ws.Range("G8").Value = "Sir William Herschel"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "1781"
$ws.Range("I8").Value = "Hydrogen, helium, and methane"

$ws.Range("F9").Value = "Neptune has the strongest winds"
$ws.Range("G9").Value = "Johann Galle and Urbain Le Verrier"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1846"
$ws.Range("I9").Value = "Hydrogen, helium, and methane"

# ---------------------------------------------------------------------------
# 3. Append the new Pluto row (row 10), matching column A's existing style
#    (bold, thin border, centered horizontally, top vertically).
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("A10").VerticalAlignment = -4160
$ws.Range("A10").Borders.LineStyle = 1

$ws.Range("B10").Value = "Pluto"
$ws.Range("C10").Value = "-229°C"
$ws.Range("D10").Value = "1,188.3 km"
$ws.Range("E10").Value = "Brownish-yellow"
$ws.Range("F10").Value = "two-thirds the width of Earth's moon"
$ws.Range("G10").Value = "Clyde Tombaugh"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "1930"
$ws.Range("I10").Value = "Mostly ice and rock"
